$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.802.52'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '2.037.36'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.31'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.605'
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.47'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.69'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '2.337.59'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.04'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.774'
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '2.026.12'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = '37.758.76'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.56'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.88'
$ws.Range('E20').Value = '  -6.71%  '
$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '223.60'
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  +2.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.41'
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '167.91'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.80'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  +8.83%  '
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.52'
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0605'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.46'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.33'
$ws.Range('E37').Value = '  +2.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.40'
$ws.Range('E38').Value = '  +3.97%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.00'
$ws.Range('E40').Value = '  +5.78%  '
$ws.Range('D41').Value = '1.535.34'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0217'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '96.17'
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0912'
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.11'
$ws.Range('E46').Value = '  -2.19%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.08'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('E48').Value = '  -0.69%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.98'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.02'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('D51').Value = '2.225.93'
$ws.Range('E51').Value = '  -1.01%  '
